$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.27
$ws.Range("B8").Value = 5.792
$ws.Range("B10").Value = 6.851999999999999
$ws.Range("B12").Value = 6.548999999999999
$ws.Range("C13").Value = -12.729
$ws.Range("B18").Value = 6.548999999999999
$ws.Range("D20").Value = -8.222
